# Scheduled-runner price refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with freshly pulled market data. A few rows also gain or lose a
# LeveProfitNQ/LeveProfitHQ cell (M/N) depending on whether NQ or HQ pricing
# is now the applicable one.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 680.1475
$ws.Range("J17").Value = 696.45764
$ws.Range("L17").Value = 2089.37292
$ws.Range("N17").Value = -2425.37292
$ws.Range("H92").Value = 661.2308
$ws.Range("I92").Value = 667.9091
$ws.Range("J92").Value = 624.5
$ws.Range("K92").Value = 667.9091
$ws.Range("L92").Value = 624.5
$ws.Range("M92").Value = 580.0909
$ws.Range("N92").Value = -3120.5
$ws.Range("H116").Value = 6360.143
$ws.Range("J116").Value = 7383
$ws.Range("L116").Value = 7383
$ws.Range("N116").Value = -14267
$ws.Range("H135").Value = 2441.6875
$ws.Range("I135").Value = 2286.2727
$ws.Range("K135").Value = 20576.4543
$ws.Range("M135").Value = -18041.4543
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 757.1429000000001
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 1150
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 1150
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -1376
$ws.Range("H32").Value = 6960.9473
$ws.Range("I32").Value = 3432.1794
$ws.Range("J32").Value = 14606.611
$ws.Range("K32").Value = 3432.1794
$ws.Range("L32").Value = 14606.611
$ws.Range("M32").Value = -3145.1794
$ws.Range("N32").Value = -15180.611
$ws.Range("H45").Value = 15627951
$ws.Range("I45").Value = 3641.8
$ws.Range("J45").Value = 41668468
$ws.Range("K45").Value = 3641.8
$ws.Range("L45").Value = 41668468
$ws.Range("M45").Value = -3264.8
$ws.Range("N45").Value = -41669222
$ws.Range("H74").Value = 2442.1052
$ws.Range("I74").Value = 1697.1666
$ws.Range("J74").Value = 3719.1428
$ws.Range("K74").Value = 1697.1666
$ws.Range("L74").Value = 3719.1428
$ws.Range("M74").Value = -823.1666
$ws.Range("N74").Value = -5467.1428
$ws.Range("H77").Value = 2442.1052
$ws.Range("I77").Value = 1697.1666
$ws.Range("J77").Value = 3719.1428
$ws.Range("K77").Value = 8485.833000000001
$ws.Range("L77").Value = 18595.714
$ws.Range("M77").Value = -4117.833000000001
$ws.Range("N77").Value = -27331.714
$ws.Range("H88").Value = 1298.9412
$ws.Range("I88").Value = 1138.8572
$ws.Range("J88").Value = 1411
$ws.Range("K88").Value = 1138.8572
$ws.Range("L88").Value = 1411
$ws.Range("M88").Value = -732.8571999999999
$ws.Range("N88").Value = -2223
$ws.Range("H91").Value = 1298.9412
$ws.Range("I91").Value = 1138.8572
$ws.Range("J91").Value = 1411
$ws.Range("K91").Value = 1138.8572
$ws.Range("L91").Value = 1411
$ws.Range("M91").Value = 265.1428000000001
$ws.Range("N91").Value = -4219
$ws.Range("H97").Value = 817.25
$ws.Range("I97").Value = 817.25
$ws.Range("K97").Value = 817.25
$ws.Range("M97").Value = -321.25
$ws.Range("H116").Value = 757.1429000000001
$ws.Range("I116").Value = 600
$ws.Range("J116").Value = 1150
$ws.Range("K116").Value = 600
$ws.Range("L116").Value = 1150
$ws.Range("M116").Value = 1694
$ws.Range("N116").Value = -5738
$ws.Range("H121").Value = 52992.4
$ws.Range("J121").Value = 52992.4
$ws.Range("L121").Value = 52992.4
$ws.Range("N121").Value = -56486.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 757.1429000000001
$ws.Range("I3").Value = 600
$ws.Range("J3").Value = 1150
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 1150
$ws.Range("M3").Value = -486
$ws.Range("N3").Value = -1378
$ws.Range("H94").Value = 2933.3333
$ws.Range("I94").Value = 2933.3333
$ws.Range("K94").Value = 2933.3333
$ws.Range("M94").Value = -2482.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7000
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 9400
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 9400
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -12894
$ws.Range("H132").Value = 1752.7916
$ws.Range("I132").Value = 1585.2858
$ws.Range("J132").Value = 1987.3
$ws.Range("K132").Value = 4755.857400000001
$ws.Range("L132").Value = 5961.9
$ws.Range("M132").Value = -2225.857400000001
$ws.Range("N132").Value = -11021.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1302
$ws.Range("I140").Value = 1015.3043
$ws.Range("K140").Value = 3045.9129
$ws.Range("M140").Value = 2134.0871
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1848.5555
$ws.Range("I102").Value = 1810.2354
$ws.Range("K102").Value = 1810.2354
$ws.Range("M102").Value = -188.2354
$ws.Range("H113").Value = 2228.5715
$ws.Range("I113").Value = 1733.3334
$ws.Range("K113").Value = 1733.3334
$ws.Range("M113").Value = 436.6666
$ws.Range("H122").Value = 1121225.2
$ws.Range("I122").Value = 1438574.2
$ws.Range("J122").Value = 10503.5
$ws.Range("K122").Value = 4315722.6
$ws.Range("L122").Value = 31510.5
$ws.Range("M122").Value = -4313272.6
$ws.Range("N122").Value = -36410.5
$ws.Range("H135").Value = 51230.75
$ws.Range("J135").Value = 51230.75
$ws.Range("L135").Value = 51230.75
$ws.Range("N135").Value = -61370.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1769.125
$ws.Range("I22").Value = 1475.25
$ws.Range("K22").Value = 1475.25
$ws.Range("M22").Value = -1180.25
$ws.Range("H27").Value = 1769.125
$ws.Range("I27").Value = 1475.25
$ws.Range("K27").Value = 1475.25
$ws.Range("M27").Value = -1368.25
$ws.Range("H46").Value = 3619
$ws.Range("I46").Value = 2561.75
$ws.Range("J46").Value = 3841.5789
$ws.Range("K46").Value = 2561.75
$ws.Range("L46").Value = 3841.5789
$ws.Range("M46").Value = -2373.75
$ws.Range("N46").Value = -4217.5789
$ws.Range("H68").Value = 6002
$ws.Range("I68").Value = 6002
$ws.Range("K68").Value = 6002
$ws.Range("M68").Value = -5253
$ws.Range("H71").Value = 6002
$ws.Range("I71").Value = 6002
$ws.Range("K71").Value = 30010
$ws.Range("M71").Value = -26266
$ws.Range("H82").Value = 2399.6667
$ws.Range("I82").Value = 2471
$ws.Range("K82").Value = 2471
$ws.Range("M82").Value = -2110
$ws.Range("H85").Value = 2399.6667
$ws.Range("I85").Value = 2471
$ws.Range("K85").Value = 2471
$ws.Range("M85").Value = -1223
$ws.Range("H93").Value = 2005.3914
$ws.Range("J93").Value = 3214.6667
$ws.Range("L93").Value = 3214.6667
$ws.Range("N93").Value = -5710.6667
$ws.Range("H94").Value = 27500
$ws.Range("J94").Value = 27500
$ws.Range("L94").Value = 27500
$ws.Range("N94").Value = -28852
$ws.Range("H100").Value = 14285.286
$ws.Range("I100").Value = 16399.8
$ws.Range("K100").Value = 16399.8
$ws.Range("M100").Value = -15858.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 62500
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -12122
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 62500
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 50000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -60608
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 12699.7
$ws.Range("I126").Value = 15642.286
$ws.Range("J126").Value = 5833.6665
$ws.Range("K126").Value = 46926.858
$ws.Range("L126").Value = 17500.9995
$ws.Range("M126").Value = -44456.858
$ws.Range("N126").Value = -22440.9995
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 1592.2727
$ws.Range("I132").Value = 1235.6875
$ws.Range("K132").Value = 3707.0625
$ws.Range("M132").Value = -1177.0625
